# Bug fixes and optimization
# Applies the "Model Performances.xlsx" update:
#  - Maps 2 RGB: fill in PSNR/SSIM data for V3.00.5-V3.00.8 (rows 7-10),
#    add 4 new model versions V3.01.5-V3.01.8 (rows 11-14) with their data,
#    add Top10% conditional formatting across columns B:I, and update the
#    active selection.
#  - GTA V Images: fill in the "Ours" row (row 8) Albedo values for Image 0
#    and add the Image 1 Albedo/Shading values, update the active selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Maps 2 RGB"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Maps 2 RGB")

# Existing rows 7-10 (V3.00.5 .. V3.00.8) were missing their PSNR/SSIM
# values - fill them in.
$ws3.Range("B7").Value = 14.69
$ws3.Range("C7").Value = 0.75209999999999999
$ws3.Range("D7").Value = 16.793600000000001
$ws3.Range("E7").Value = 0.81799999999999995
$ws3.Range("F7").Value = 16.810199999999998
$ws3.Range("G7").Value = 0.83389999999999997
$ws3.Range("H7").Value = 22.0181
$ws3.Range("I7").Value = 0.84430000000000005

$ws3.Range("B8").Value = 11.166499999999999
$ws3.Range("C8").Value = 0.69379999999999997
$ws3.Range("D8").Value = 15.355600000000001
$ws3.Range("E8").Value = 0.80510000000000004
$ws3.Range("F8").Value = 15.696899999999999
$ws3.Range("G8").Value = 0.79110000000000003
$ws3.Range("H8").Value = 15.081099999999999
$ws3.Range("I8").Value = 0.77100000000000002

$ws3.Range("B9").Value = 14.463699999999999
$ws3.Range("C9").Value = 0.75170000000000003
$ws3.Range("D9").Value = 16.669699999999999
$ws3.Range("E9").Value = 0.83250000000000002
$ws3.Range("F9").Value = 16.628599999999999
$ws3.Range("G9").Value = 0.8004
$ws3.Range("H9").Value = 18.977599999999999
$ws3.Range("I9").Value = 0.83520000000000005

$ws3.Range("B10").Value = 12.6652
$ws3.Range("C10").Value = 0.73570000000000002
$ws3.Range("D10").Value = 15.749000000000001
$ws3.Range("E10").Value = 0.83179999999999998
$ws3.Range("F10").Value = 13.571300000000001
$ws3.Range("G10").Value = 0.78620000000000001
$ws3.Range("H10").Value = 16.597100000000001
$ws3.Range("I10").Value = 0.80189999999999995

# Apply the "0.0000" number format (matches the style already used by
# the data cells elsewhere in the workbook) to the newly-populated cells.
$ws3.Range("B7:I10").NumberFormat = "0.0000"

# New model versions V3.01.5 .. V3.01.8 (rows 11-14)
$ws3.Range("A11").Value = "V3.01.5"
$ws3.Range("B11").Value = 14.2508
$ws3.Range("C11").Value = 0.7399
$ws3.Range("D11").Value = 12.297000000000001
$ws3.Range("E11").Value = 0.73119999999999996
$ws3.Range("F11").Value = 12.095800000000001
$ws3.Range("G11").Value = 0.74939999999999996
$ws3.Range("H11").Value = 19.1495
$ws3.Range("I11").Value = 0.8

$ws3.Range("A12").Value = "V3.01.6"
$ws3.Range("B12").Value = 14.7559
$ws3.Range("C12").Value = 0.77680000000000005
$ws3.Range("D12").Value = 13.9665
$ws3.Range("E12").Value = 0.81599999999999995
$ws3.Range("F12").Value = 17.352699999999999
$ws3.Range("G12").Value = 0.82830000000000004
$ws3.Range("H12").Value = 24.509899999999998
$ws3.Range("I12").Value = 0.88629999999999998

$ws3.Range("A13").Value = "V3.01.7"
$ws3.Range("B13").Value = 13.792299999999999
$ws3.Range("C13").Value = 0.77380000000000004
$ws3.Range("D13").Value = 13.6998
$ws3.Range("E13").Value = 0.83379999999999999
$ws3.Range("F13").Value = 16.271899999999999
$ws3.Range("G13").Value = 0.83140000000000003
$ws3.Range("H13").Value = 13.792299999999999
$ws3.Range("I13").Value = 0.89549999999999996

$ws3.Range("A14").Value = "V3.01.8"
$ws3.Range("B14").Value = 14.797700000000001
$ws3.Range("C14").Value = 0.79020000000000001
$ws3.Range("D14").Value = 16.731000000000002
$ws3.Range("E14").Value = 0.85550000000000004
$ws3.Range("F14").Value = 11.760999999999999
$ws3.Range("G14").Value = 0.79979999999999996
$ws3.Range("H14").Value = 19.196899999999999
$ws3.Range("I14").Value = 0.8659

$ws3.Range("B11:I14").NumberFormat = "0.0000"

# New "Top 10%" conditional-formatting rules across the full B:I columns
# (adds on top of the existing per-header-cell rules).
$cols = @("B","C","D","E","F","G","H","I")
$newConds = @()
foreach ($col in $cols) {
    $rng = $ws3.Range($col + "1:" + $col + "1048576")
    $cond = $rng.FormatConditions.AddTop10()
    $cond.Percent = $true
    $cond.Rank = 10
    $cond.Font.Color = 24832
    $cond.Interior.Color = 13561798
    $newConds += $cond
}
$newPriorities = @(8,7,6,5,4,3,2,1)
for ($i = 0; $i -lt $newConds.Length; $i++) {
    $newConds[$i].Priority = $newPriorities[$i]
}

# Existing header-cell rules get pushed back (lower priority) now that the
# new column-wide rules take precedence.
$oldCells = @("D2","E2","F2","G2","B2","C2","H2","I2")
$oldPriorities = @(16,15,14,13,12,11,10,9)
for ($i = 0; $i -lt $oldCells.Length; $i++) {
    $fc = $ws3.Range($oldCells[$i]).FormatConditions.Item(1)
    $fc.Priority = $oldPriorities[$i]
}

# Update selection to match the edited area.
$ws3.Activate()
$ws3.Range("H12").Select()

# ---------------------------------------------------------------------
# Sheet "GTA V Images"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("GTA V Images")

# Row 8 ("Ours"): fill in Image 0 Albedo values and add Image 1 values.
$ws4.Range("B8").Value = 8.2411999999999992
$ws4.Range("C8").Value = 0.64659999999999995
$ws4.Range("G8").Value = 8.343
$ws4.Range("H8").Value = 0.70499999999999996
$ws4.Range("I8").Value = 8.1175999999999995
$ws4.Range("J8").Value = 0.70920000000000005
$ws4.Range("B8:C8,G8:J8").NumberFormat = "0.0000"

# Update selection / view to match.
$ws4.Activate()
$ws4.Range("O14").Select()

Write-Output "edit.ps1 completed"
